$wb = $excel.ActiveWorkbook

# "normal" sheet: add new column S = INT_recruitTime
$ws1 = $wb.Worksheets.Item("normal")
$ws1.Range("S1").Value = "INT_recruitTime"

$ws1.Range("S2").Value = 30
$ws1.Range("S3").Value = 30
$ws1.Range("S4").Value = 30
$ws1.Range("S5").Value = 35
$ws1.Range("S6").Value = 35
$ws1.Range("S7").Value = 35
$ws1.Range("S8").Value = 50
$ws1.Range("S9").Value = 50
$ws1.Range("S10").Value = 50
$ws1.Range("S11").Value = 55
$ws1.Range("S12").Value = 55
$ws1.Range("S13").Value = 55
$ws1.Range("S14").Value = 110
$ws1.Range("S15").Value = 110
$ws1.Range("S16").Value = 110
$ws1.Range("S17").Value = 120
$ws1.Range("S18").Value = 120
$ws1.Range("S19").Value = 120
$ws1.Range("S20").Value = 180
$ws1.Range("S21").Value = 180
$ws1.Range("S22").Value = 180
$ws1.Range("S23").Value = 190
$ws1.Range("S24").Value = 190
$ws1.Range("S25").Value = 190

# "special" sheet: add new column P = INT_recruitTime
$ws2 = $wb.Worksheets.Item("special")
$ws2.Range("P1").Value = "INT_recruitTime"

$ws2.Range("P2").Value = 40
$ws2.Range("P3").Value = 60
$ws2.Range("P4").Value = 80
$ws2.Range("P5").Value = 95
$ws2.Range("P6").Value = 140
$ws2.Range("P7").Value = 180
$ws2.Range("P8").Value = 240
$ws2.Range("P9").Value = 280
